$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column A into column B for every used row (1-136)
$ws.Range("A1:A136").Copy($ws.Range("B1:B136"))

# Remove the old stray cells that used to live in columns C and D
# (row 6 col C, rows 39/42/44/51/56/58 col D) so the used range
# shrinks back down to just A:B.
$ws.Range("C1:D136").Clear()
